# Auto-generated script applying the week10 spreads_tracker update.
# Rebuilds rows 2:53 (A:J) with the new sort order / new games, and
# extends the sheet from 45 to 53 data rows, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 52,10
$data[0,0] = 0
$data[0,1] = 'Utah'
$data[0,2] = 'Cincinnati'
$data[0,3] = 9.2
$data[0,4] = 5.800000000000001
$data[0,5] = 'Utah -5.5'
$data[0,6] = 'Utah -6.5'
$data[0,7] = 'Utah -12.3'
$data[0,8] = 12.3
$data[0,9] = 6.5
$data[1,0] = 1
$data[1,1] = 'Louisiana Tech'
$data[1,2] = 'Sam Houston'
$data[1,3] = 2.5
$data[1,4] = 5.800000000000001
$data[1,5] = 'Louisiana Tech -17.5'
$data[1,6] = 'Louisiana Tech -16.5'
$data[1,7] = 'Louisiana Tech -22.3'
$data[1,8] = 22.3
$data[1,9] = 16.5
$data[2,0] = 2
$data[2,1] = 'Syracuse'
$data[2,2] = 'North Carolina'
$data[2,3] = 6.9
$data[2,4] = 5.7
$data[2,5] = 'Syracuse -nan'
$data[2,6] = 'Syracuse -1.5'
$data[2,7] = 'Syracuse -7.2'
$data[2,8] = 7.2
$data[2,9] = 1.5
$data[3,0] = 3
$data[3,1] = 'Rice'
$data[3,2] = 'Memphis'
$data[3,3] = 3.5
$data[3,4] = 5.600000000000001
$data[3,5] = 'Memphis -17.5'
$data[3,6] = 'Memphis -14.5'
$data[3,7] = 'Memphis -20.1'
$data[3,8] = -20.1
$data[3,9] = -14.5
$data[4,0] = 4
$data[4,1] = 'UConn'
$data[4,2] = 'UAB'
$data[4,3] = 4.6
$data[4,4] = 5
$data[4,5] = 'UConn -12.5'
$data[4,6] = 'UConn -11.5'
$data[4,7] = 'UConn -16.5'
$data[4,8] = 16.5
$data[4,9] = 11.5
$data[5,0] = 5
$data[5,1] = 'California'
$data[5,2] = 'Virginia'
$data[5,3] = 6.9
$data[5,4] = 4.4
$data[5,5] = 'Virginia -4.0'
$data[5,6] = 'Virginia -4'
$data[5,7] = 'Virginia -8.4'
$data[5,8] = -8.4
$data[5,9] = -4
$data[6,0] = 6
$data[6,1] = 'Ohio State'
$data[6,2] = 'Penn State'
$data[6,3] = 9.3
$data[6,4] = 4.199999999999999
$data[6,5] = 'Ohio State -4.0'
$data[6,6] = 'Ohio State -20.5'
$data[6,7] = 'Ohio State -16.3'
$data[6,8] = 16.3
$data[6,9] = 20.5
$data[7,0] = 7
$data[7,1] = 'Florida State'
$data[7,2] = 'Wake Forest'
$data[7,3] = 8.1
$data[7,4] = 4
$data[7,5] = 'Florida State -7.0'
$data[7,6] = 'Florida State -8.5'
$data[7,7] = 'Florida State -12.5'
$data[7,8] = 12.5
$data[7,9] = 8.5
$data[8,0] = 8
$data[8,1] = 'Baylor'
$data[8,2] = 'UCF'
$data[8,3] = 8.9
$data[8,4] = 4
$data[8,5] = 'Baylor -4.5'
$data[8,6] = 'Baylor -5.5'
$data[8,7] = 'Baylor -1.5'
$data[8,8] = 1.5
$data[8,9] = 5.5
$data[9,0] = 9
$data[9,1] = 'Auburn'
$data[9,2] = 'Kentucky'
$data[9,3] = 8.1
$data[9,4] = 3.9
$data[9,5] = 'Auburn -10.0'
$data[9,6] = 'Auburn -10'
$data[9,7] = 'Auburn -13.9'
$data[9,8] = 13.9
$data[9,9] = 10
$data[10,0] = 10
$data[10,1] = 'Kennesaw State'
$data[10,2] = 'UTEP'
$data[10,3] = 4.8
$data[10,4] = 3.800000000000001
$data[10,5] = 'Kennesaw State -10.0'
$data[10,6] = 'Kennesaw State -8.5'
$data[10,7] = 'Kennesaw State -12.3'
$data[10,8] = 12.3
$data[10,9] = 8.5
$data[11,0] = 11
$data[11,1] = 'Texas'
$data[11,2] = 'Vanderbilt'
$data[11,3] = 9.7
$data[11,4] = 3.7
$data[11,5] = 'Texas -3.5'
$data[11,6] = 'Texas -2.5'
$data[11,7] = 'Texas -6.2'
$data[11,8] = 6.2
$data[11,9] = 2.5
$data[12,0] = 12
$data[12,1] = 'Troy'
$data[12,2] = 'Arkansas State'
$data[12,3] = 5.6
$data[12,4] = 3.6
$data[12,5] = 'Troy -7.0'
$data[12,6] = 'Troy -7.5'
$data[12,7] = 'Troy -11.1'
$data[12,8] = 11.1
$data[12,9] = 7.5
$data[13,0] = 13
$data[13,1] = 'Florida'
$data[13,2] = 'Georgia'
$data[13,3] = 9.1
$data[13,4] = 3.4
$data[13,5] = 'Georgia -7.5'
$data[13,6] = 'Georgia -7.5'
$data[13,7] = 'Georgia -10.9'
$data[13,8] = -10.9
$data[13,9] = -7.5
$data[14,0] = 14
$data[14,1] = 'UL Monroe'
$data[14,2] = 'Old Dominion'
$data[14,3] = 3
$data[14,4] = 3.199999999999999
$data[14,5] = 'Old Dominion -14.0'
$data[14,6] = 'Old Dominion -14'
$data[14,7] = 'Old Dominion -17.2'
$data[14,8] = -17.2
$data[14,9] = -14
$data[15,0] = 15
$data[15,1] = 'Michigan'
$data[15,2] = 'Purdue'
$data[15,3] = 6
$data[15,4] = 3.100000000000001
$data[15,5] = 'Michigan -20.5'
$data[15,6] = 'Michigan -21'
$data[15,7] = 'Michigan -24.1'
$data[15,8] = 24.1
$data[15,9] = 21
$data[16,0] = 16
$data[16,1] = 'North Texas'
$data[16,2] = 'Navy'
$data[16,3] = 8.3
$data[16,4] = 2.9
$data[16,5] = 'North Texas -4.5'
$data[16,6] = 'North Texas -6.5'
$data[16,7] = 'North Texas -9.4'
$data[16,8] = 9.4
$data[16,9] = 6.5
$data[17,0] = 17
$data[17,1] = 'Coastal Carolina'
$data[17,2] = 'Marshall'
$data[17,3] = 4.9
$data[17,4] = 2.9
$data[17,5] = 'Marshall -3.0'
$data[17,6] = 'Marshall -4'
$data[17,7] = 'Marshall -6.9'
$data[17,8] = -6.9
$data[17,9] = -4
$data[18,0] = 18
$data[18,1] = 'Boston College'
$data[18,2] = 'Notre Dame'
$data[18,3] = 4.1
$data[18,4] = 2.800000000000001
$data[18,5] = 'Notre Dame -27.5'
$data[18,6] = 'Notre Dame -27.5'
$data[18,7] = 'Notre Dame -30.3'
$data[18,8] = -30.3
$data[18,9] = -27.5
$data[19,0] = 19
$data[19,1] = 'Tennessee'
$data[19,2] = 'Oklahoma'
$data[19,3] = 9.7
$data[19,4] = 2.8
$data[19,5] = 'Tennessee -2.5'
$data[19,6] = 'Tennessee -4'
$data[19,7] = 'Tennessee -1.2'
$data[19,8] = 1.2
$data[19,9] = 4
$data[20,0] = 20
$data[20,1] = 'Iowa State'
$data[20,2] = 'Arizona State'
$data[20,3] = 9
$data[20,4] = 2.5
$data[20,5] = 'Iowa State -3.0'
$data[20,6] = 'Iowa State -3.5'
$data[20,7] = 'Iowa State -6.0'
$data[20,8] = 6
$data[20,9] = 3.5
$data[21,0] = 21
$data[21,1] = 'Stanford'
$data[21,2] = 'Pittsburgh'
$data[21,3] = 5.1
$data[21,4] = 2.5
$data[21,5] = 'Pittsburgh -15.5'
$data[21,6] = 'Pittsburgh -14.5'
$data[21,7] = 'Pittsburgh -17.0'
$data[21,8] = -17
$data[21,9] = -14.5
$data[22,0] = 22
$data[22,1] = 'Missouri State'
$data[22,2] = 'Florida International'
$data[22,3] = 5.1
$data[22,4] = 2.4
$data[22,5] = 'Missouri State -2.5'
$data[22,6] = 'Missouri State -3.5'
$data[22,7] = 'Missouri State -5.9'
$data[22,8] = 5.9
$data[22,9] = 3.5
$data[23,0] = 23
$data[23,1] = 'Kansas State'
$data[23,2] = 'Texas Tech'
$data[23,3] = 8.9
$data[23,4] = 2.300000000000001
$data[23,5] = 'Kansas State -4.5'
$data[23,6] = 'Texas Tech -7'
$data[23,7] = 'Texas Tech -9.3'
$data[23,8] = -9.3
$data[23,9] = -7
$data[24,0] = 24
$data[24,1] = 'Oregon State'
$data[24,2] = 'Washington State'
$data[24,3] = 5.9
$data[24,4] = 2.3
$data[24,5] = 'Washington State -3.0'
$data[24,6] = 'Washington State -3.5'
$data[24,7] = 'Washington State -5.8'
$data[24,8] = -5.8
$data[24,9] = -3.5
$data[25,0] = 25
$data[25,1] = 'Boise State'
$data[25,2] = 'Fresno State'
$data[25,3] = 5.4
$data[25,4] = 2.199999999999999
$data[25,5] = 'Boise State -17.5'
$data[25,6] = 'Boise State -17.5'
$data[25,7] = 'Boise State -19.7'
$data[25,8] = 19.7
$data[25,9] = 17.5
$data[26,0] = 26
$data[26,1] = 'Maryland'
$data[26,2] = 'Indiana'
$data[26,3] = 7.5
$data[26,4] = 2.199999999999999
$data[26,5] = 'Indiana -16.5'
$data[26,6] = 'Indiana -22.5'
$data[26,7] = 'Indiana -20.3'
$data[26,8] = -20.3
$data[26,9] = -22.5
$data[27,0] = 27
$data[27,1] = 'Western Kentucky'
$data[27,2] = 'New Mexico State'
$data[27,3] = 4.5
$data[27,4] = 2.1
$data[27,5] = 'Western Kentucky -9.5'
$data[27,6] = 'Western Kentucky -9.5'
$data[27,7] = 'Western Kentucky -11.6'
$data[27,8] = 11.6
$data[27,9] = 9.5
$data[28,0] = 28
$data[28,1] = 'Virginia Tech'
$data[28,2] = 'Louisville'
$data[28,3] = 6.8
$data[28,4] = 2.1
$data[28,5] = 'Louisville -11.5'
$data[28,6] = 'Louisville -10.5'
$data[28,7] = 'Louisville -12.6'
$data[28,8] = -12.6
$data[28,9] = -10.5
$data[29,0] = 29
$data[29,1] = 'Colorado'
$data[29,2] = 'Arizona'
$data[29,3] = 8.3
$data[29,4] = 2
$data[29,5] = 'Arizona -4.5'
$data[29,6] = 'Arizona -4.5'
$data[29,7] = 'Arizona -2.5'
$data[29,8] = -2.5
$data[29,9] = -4.5
$data[30,0] = 30
$data[30,1] = 'San Diego State'
$data[30,2] = 'Wyoming'
$data[30,3] = 6.4
$data[30,4] = 1.9
$data[30,5] = 'San Diego State -10.0'
$data[30,6] = 'San Diego State -11.5'
$data[30,7] = 'San Diego State -13.4'
$data[30,8] = 13.4
$data[30,9] = 11.5
$data[31,0] = 31
$data[31,1] = 'Middle Tennessee'
$data[31,2] = 'Jacksonville State'
$data[31,3] = 3.8
$data[31,4] = 1.6
$data[31,5] = 'Jacksonville State -6.5'
$data[31,6] = 'Jacksonville State -6'
$data[31,7] = 'Jacksonville State -7.6'
$data[31,8] = -7.6
$data[31,9] = -6
$data[32,0] = 32
$data[32,1] = 'Ole Miss'
$data[32,2] = 'South Carolina'
$data[32,3] = 8.6
$data[32,4] = 1.5
$data[32,5] = 'Ole Miss -13.5'
$data[32,6] = 'Ole Miss -13.5'
$data[32,7] = 'Ole Miss -15.0'
$data[32,8] = 15
$data[32,9] = 13.5
$data[33,0] = 33
$data[33,1] = 'Minnesota'
$data[33,2] = 'Michigan State'
$data[33,3] = 8.1
$data[33,4] = 1.5
$data[33,5] = 'Minnesota -5.5'
$data[33,6] = 'Minnesota -5.5'
$data[33,7] = 'Minnesota -7.0'
$data[33,8] = 7
$data[33,9] = 5.5
$data[34,0] = 34
$data[34,1] = 'UNLV'
$data[34,2] = 'New Mexico'
$data[34,3] = 7.6
$data[34,4] = 1.5
$data[34,5] = 'UNLV -5.5'
$data[34,6] = 'UNLV -4.5'
$data[34,7] = 'UNLV -6.0'
$data[34,8] = 6
$data[34,9] = 4.5
$data[35,0] = 35
$data[35,1] = 'UTSA'
$data[35,2] = 'Tulane'
$data[35,3] = 7.7
$data[35,4] = 1.5
$data[35,5] = 'Tulane -3.5'
$data[35,6] = 'Tulane -4.5'
$data[35,7] = 'Tulane -3.0'
$data[35,8] = -3
$data[35,9] = -4.5
$data[36,0] = 36
$data[36,1] = 'SMU'
$data[36,2] = 'Miami'
$data[36,3] = 8.8
$data[36,4] = 1.4
$data[36,5] = 'Miami -1.5'
$data[36,6] = 'Miami -10.5'
$data[36,7] = 'Miami -9.1'
$data[36,8] = -9.1
$data[36,9] = -10.5
$data[37,0] = 37
$data[37,1] = 'Nebraska'
$data[37,2] = 'USC'
$data[37,3] = 9.3
$data[37,4] = 1.4
$data[37,5] = 'USC -4.5'
$data[37,6] = 'USC -6'
$data[37,7] = 'USC -4.6'
$data[37,8] = -4.6
$data[37,9] = -6
$data[38,0] = 38
$data[38,1] = 'NC State'
$data[38,2] = 'Georgia Tech'
$data[38,3] = 8.1
$data[38,4] = 1.2
$data[38,5] = 'Georgia Tech -4.5'
$data[38,6] = 'Georgia Tech -6.5'
$data[38,7] = 'Georgia Tech -7.7'
$data[38,8] = -7.7
$data[38,9] = -6.5
$data[39,0] = 39
$data[39,1] = 'Liberty'
$data[39,2] = 'Delaware'
$data[39,3] = 6.3
$data[39,4] = 1.1
$data[39,5] = 'Liberty -2.5'
$data[39,6] = 'Liberty -3'
$data[39,7] = 'Liberty -4.1'
$data[39,8] = 4.1
$data[39,9] = 3
$data[40,0] = 40
$data[40,1] = 'Kansas'
$data[40,2] = 'Oklahoma State'
$data[40,3] = 3.2
$data[40,4] = 1
$data[40,5] = 'Kansas -24.5'
$data[40,6] = 'Kansas -25.5'
$data[40,7] = 'Kansas -26.5'
$data[40,8] = 26.5
$data[40,9] = 25.5
$data[41,0] = 41
$data[41,1] = 'Clemson'
$data[41,2] = 'Duke'
$data[41,3] = 9.2
$data[41,4] = 1
$data[41,5] = 'Clemson -3.5'
$data[41,6] = 'Clemson -3.5'
$data[41,7] = 'Clemson -4.5'
$data[41,8] = 4.5
$data[41,9] = 3.5
$data[42,0] = 42
$data[42,1] = 'Arkansas'
$data[42,2] = 'Mississippi State'
$data[42,3] = 9.3
$data[42,4] = 0.9000000000000004
$data[42,5] = 'Arkansas -3.5'
$data[42,6] = 'Arkansas -3.5'
$data[42,7] = 'Arkansas -4.4'
$data[42,8] = 4.4
$data[42,9] = 3.5
$data[43,0] = 43
$data[43,1] = 'South Alabama'
$data[43,2] = 'Louisiana'
$data[43,3] = 6.2
$data[43,4] = 0.7999999999999998
$data[43,5] = 'South Alabama -4.5'
$data[43,6] = 'South Alabama -4'
$data[43,7] = 'South Alabama -4.8'
$data[43,8] = 4.8
$data[43,9] = 4
$data[44,0] = 44
$data[44,1] = 'Houston'
$data[44,2] = 'West Virginia'
$data[44,3] = 6.5
$data[44,4] = 0.4000000000000004
$data[44,5] = 'Houston -15.5'
$data[44,6] = 'Houston -15.5'
$data[44,7] = 'Houston -15.1'
$data[44,8] = 15.1
$data[44,9] = 15.5
$data[45,0] = 45
$data[45,1] = 'Illinois'
$data[45,2] = 'Rutgers'
$data[45,3] = 8.1
$data[45,4] = 0.3000000000000007
$data[45,5] = 'Illinois -10.5'
$data[45,6] = 'Illinois -12.5'
$data[45,7] = 'Illinois -12.8'
$data[45,8] = 12.8
$data[45,9] = 12.5
$data[46,0] = 46
$data[46,1] = 'Texas State'
$data[46,2] = 'James Madison'
$data[46,3] = 7
$data[46,4] = 0.2000000000000002
$data[46,5] = 'James Madison -6.5'
$data[46,6] = 'James Madison -6.5'
$data[46,7] = 'James Madison -6.7'
$data[46,8] = -6.7
$data[46,9] = -6.5
$data[47,0] = 47
$data[47,1] = 'Temple'
$data[47,2] = 'East Carolina'
$data[47,3] = 7.2
$data[47,4] = 0.2000000000000002
$data[47,5] = 'East Carolina -4.5'
$data[47,6] = 'East Carolina -4.5'
$data[47,7] = 'East Carolina -4.7'
$data[47,8] = -4.7
$data[47,9] = -4.5
$data[48,0] = 48
$data[48,1] = 'Bowling Green'
$data[48,2] = 'Buffalo'
$data[48,3] = 6.1
$data[48,4] = 0.2
$data[48,5] = 'Bowling Green -1.5'
$data[48,6] = 'Bowling Green -1.5'
$data[48,7] = 'Bowling Green -1.7'
$data[48,8] = 1.7
$data[48,9] = 1.5
$data[49,0] = 49
$data[49,1] = 'Western Michigan'
$data[49,2] = 'Central Michigan'
$data[49,3] = 6.3
$data[49,4] = 0
$data[49,5] = 'Western Michigan -6.0'
$data[49,6] = 'Western Michigan -6.5'
$data[49,7] = 'Western Michigan -6.5'
$data[49,8] = 6.5
$data[49,9] = 6.5
$data[50,0] = 50
$data[50,1] = 'San José State'
$data[50,2] = 'Hawai''i'
$data[50,3] = 7
$data[50,4] = 0
$data[50,5] = 'San José State -1.5'
$data[50,6] = 'San José State -1.5'
$data[50,7] = 'San José State -1.5'
$data[50,8] = 1.5
$data[50,9] = 1.5
$data[51,0] = 51
$data[51,1] = 'Air Force'
$data[51,2] = 'Army'
$data[51,3] = 7
$data[51,4] = 0
$data[51,5] = 'Army -1.5'
$data[51,6] = 'Army -1.5'
$data[51,7] = 'Army -1.5'
$data[51,8] = -1.5
$data[51,9] = -1.5

$ws.Range("A2:J53").Value = $data

# Column P ("PEAR_OU") is always 0 for every data row.
$pdata = New-Object 'object[,]' 52,1
$pdata[0,0] = 0
$pdata[1,0] = 0
$pdata[2,0] = 0
$pdata[3,0] = 0
$pdata[4,0] = 0
$pdata[5,0] = 0
$pdata[6,0] = 0
$pdata[7,0] = 0
$pdata[8,0] = 0
$pdata[9,0] = 0
$pdata[10,0] = 0
$pdata[11,0] = 0
$pdata[12,0] = 0
$pdata[13,0] = 0
$pdata[14,0] = 0
$pdata[15,0] = 0
$pdata[16,0] = 0
$pdata[17,0] = 0
$pdata[18,0] = 0
$pdata[19,0] = 0
$pdata[20,0] = 0
$pdata[21,0] = 0
$pdata[22,0] = 0
$pdata[23,0] = 0
$pdata[24,0] = 0
$pdata[25,0] = 0
$pdata[26,0] = 0
$pdata[27,0] = 0
$pdata[28,0] = 0
$pdata[29,0] = 0
$pdata[30,0] = 0
$pdata[31,0] = 0
$pdata[32,0] = 0
$pdata[33,0] = 0
$pdata[34,0] = 0
$pdata[35,0] = 0
$pdata[36,0] = 0
$pdata[37,0] = 0
$pdata[38,0] = 0
$pdata[39,0] = 0
$pdata[40,0] = 0
$pdata[41,0] = 0
$pdata[42,0] = 0
$pdata[43,0] = 0
$pdata[44,0] = 0
$pdata[45,0] = 0
$pdata[46,0] = 0
$pdata[47,0] = 0
$pdata[48,0] = 0
$pdata[49,0] = 0
$pdata[50,0] = 0
$pdata[51,0] = 0
$ws.Range("P2:P53").Value = $pdata

# Extend column A formatting (bold, centered, thin border) down to the
# newly added rows 46:53, matching the style already used in A2:A45.
$ws.Range("A45").Copy() | Out-Null
$ws.Range("A46:A53").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Select()
